# Weekly update: two new price records were inserted into the data table.
# The first is inserted as the new row 316 (pushing the former rows
# 316..408 down to 317..409); the second is inserted right after, as the
# new row 408 (pushing the former row 408 - now at 408 after the first
# shift - down to 410).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MangoRow($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = 4
    $ws.Cells.Item($Row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($Row, 3).Value = "Los Lagos"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 10
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100108
    $ws.Cells.Item($Row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($Row, 9).Value = 100108002
    $ws.Cells.Item($Row, 10).Value = "Mango"
    $ws.Cells.Item($Row, 11).Value = "Sin especificar"
    $ws.Cells.Item($Row, 12).Value = "Primera"
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = 4
}

# Insert the first new record at row 316 (date 45120, 2023-07-13).
$ws.Rows.Item(316).Insert()
Set-MangoRow 316 45120 120 8500 9000 8750 "Perú" 2188

# Insert the second new record at row 408 (date 45121, 2023-07-14).
$ws.Rows.Item(408).Insert()
Set-MangoRow 408 45121 120 8500 9000 8750 "Brasil" 2188
